# Update the two-digit x two-digit multiplication problems in the
# worksheet table. Each of the 5 "problem" rows (1, 5, 10, 15, 20) has
# 5 cells; replace the expression text in each cell in place so the
# existing run formatting (font/size) on the w:t text is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    "1,1"  = "18×91="
    "1,2"  = "62×70="
    "1,3"  = "23×23="
    "1,4"  = "28×91="
    "1,5"  = "28×91="
    "5,1"  = "98×90="
    "5,2"  = "79×83="
    "5,3"  = "98×98="
    "5,4"  = "53×12="
    "5,5"  = "62×48="
    "10,1" = "43×73="
    "10,2" = "57×27="
    "10,3" = "24×70="
    "10,4" = "86×86="
    "10,5" = "42×41="
    "15,1" = "37×12="
    "15,2" = "43×49="
    "15,3" = "46×30="
    "15,4" = "75×28="
    "15,5" = "17×11="
    "20,1" = "34×47="
    "20,2" = "47×51="
    "20,3" = "94×36="
    "20,4" = "96×53="
    "20,5" = "83×85="
}

$rowIndexes = @(1, 5, 10, 15, 20)
foreach ($r in $rowIndexes) {
    for ($c = 1; $c -le 5; $c++) {
        $key = "$r,$c"
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$key]
    }
}
